$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '52.068.59'
$ws.Range("E2").Value = '  +4.85%  '

$ws.Range("D3").Value = '2.783.53'
$ws.Range("E3").Value = '  +5.10%  '

$ws.Range("E4").Value = '  -0.01%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '342.76'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.53%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '115.64'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.64%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.548'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.69%  '

$ws.Range("E8").Value = '  -0.04%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.578'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.88%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.05'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +5.75%  '

$ws.Range("E11").Value = '  +4.90%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '20.02'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.17%  '

$ws.Range("E13").Value = '  +1.49%  '

$ws.Range("E14").Value = '  +0.34%  '

$ws.Range("D15").Value = '3.218.01'
$ws.Range("E15").Value = '  +4.97%  '

$ws.Range("D16").Value = '2.774.00'
$ws.Range("E16").Value = '  +4.95%  '

$ws.Range("D17").Value = '51.931.52'
$ws.Range("E17").Value = '  +4.57%  '

$ws.Range("E18").Value = '  +1.92%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.20'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +10.15%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.03'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.01%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.28'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.43%  '

$ws.Range("D22").Value = '0.0₃0978'
$ws.Range("E22").Value = '  +3.04%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '277.04'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.02%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.07'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.41%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.77'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +7.73%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '26.68'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.99%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.01%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.18'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.13%  '

$ws.Range("E29").Value = '  +0.71%  '

$ws.Range("E30").Value = '  +2.16%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '34.71'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.47%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '50.16'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.07%  '

$ws.Range("E33").Value = '  +4.49%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0818'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.40%  '

$ws.Range("E35").Value = '  -0.15%  '

$ws.Range("E36").Value = '  +3.37%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '18.97'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.23%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.96'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.27%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.22'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.39%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0383'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +11.13%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.68'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +26.81%  '

$ws.Range("B42").Value = 'Stellar'
$ws.Range("C42").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.116'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.57%  '

$ws.Range("B43").Value = 'WEMIXToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.34'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.52%  '

$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '23.28'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.91%  '

$ws.Range("B45").Value = 'Monero'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '125.69'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.53%  '

$ws.Range("D46").Value = '2.071.56'
$ws.Range("E46").Value = '  +0.28%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.32'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.02%  '

$ws.Range("E48").Value = '  +0.10%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.56'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.45%  '

$ws.Range("B50").Value = 'SEI'
$ws.Range("C50").Value = 'https://coinranking.com/coin/8nxCqs-uj+sei-sei'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.894'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +16.55%  '

$ws.Range("B51").Value = 'FraxShare'
$ws.Range("C51").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.88'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.75%  '
